$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B42").Value = "/bodyset/trapezoid"

$ws.Range("C42").Value = "''0 0 0"
$ws.Range("C42").Style = "Normal"

$ws.Range("C43").Value = "''-0.0047000000000000002 0.033500000000000002 -0.0089999999999999993"
$ws.Range("C43").Style = "Normal"

$ws.Range("C46").Value = "''0 0 0"
$ws.Range("C46").Style = "Normal"

$ws.Range("C47").Value = "''0.012999999999999999 0.040000000000000001 -0.0050000000000000001"
$ws.Range("C47").Style = "Normal"

$ws.Range("C48").Value = "''0 0 0"
$ws.Range("C48").Style = "Normal"

$ws.Range("C49").Value = "''0.02 0.035999999999999997 0.0050000000000000001"
$ws.Range("C49").Style = "Normal"
